# "hours update source add"
# Add two new log entries (rows 45 and 46) to the status-report sheet and
# move the active-cell selection down to the next blank row (A47).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 45: new "Group Meeting" entry, 2010-03-23 (Excel serial 40260), 2.5 hrs
$ws.Range("A45").Value = 40260
$ws.Range("B45").Value = 2.5
$ws.Range("C45").Value = "Group Meeting"

# Row 46: new "Weekly Meeting" entry, 2010-03-23 (Excel serial 40260), 1 hr
$ws.Range("A46").Value = 40260
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "Weekly Meeting"

# Move selection to the next empty row, same as the author's workbook state
$ws.Range("A47").Select()
